$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix rows 326:336 - column D (bsecode) was stored as text, convert to
# --- a real number (same digits, new type). ---
$bsecodeFix = @{
    326 = 532466
    327 = 532540
    328 = 542650
    329 = 500209
    330 = 532689
    331 = 532454
    332 = 500043
    333 = 511243
    334 = 500850
    335 = 530005
    336 = 532461
}
foreach ($r in $bsecodeFix.Keys) {
    $ws.Cells.Item($r, 4).Value = $bsecodeFix[$r]
}

# --- Append new rows 337:350 (break out stock.yaml completed) ---
$newRows = @(
    @(1,  "INDIGO",     "Interglobe Aviation Limited",                  "539448", -0.57, 4227.4,  508342,    "day", "13/08/2024 11:34:40"),
    @(2,  "INDIAMART",  "Indiamart Intermesh Ltd",                      "542726", -2.1,  2699.4,  128834,    "day", "13/08/2024 11:34:40"),
    @(3,  "NESTLEIND",  "Nestle India Limited",                         "500790", 0.47,  2484.7,  311846,    "day", "13/08/2024 11:34:40"),
    @(4,  "BAJAJFINSV", "Bajaj Finserv Limited",                        "532978", -1,    1543.35, 1817066,   "day", "13/08/2024 11:34:40"),
    @(5,  "BATAINDIA",  "Bata India Limited",                           "500043", -1.13, 1402.8,  354502,    "day", "13/08/2024 11:34:40"),
    @(6,  "TATAMOTORS", "Tata Motors Limited",                          "500570", -2.11, 1053.45, 7508231,   "day", "13/08/2024 11:34:40"),
    @(7,  "MARICO",     "Marico Limited",                               "531642", 2.47,  660.55,  4531176,   "day", "13/08/2024 11:34:40"),
    @(8,  "GUJGASLTD",  "Gujarat Gas Limited",                          "539336", -3.74, 595.75,  660459,    "day", "13/08/2024 11:34:40"),
    @(9,  "RECLTD",     "Rural Electrification Corporation Limited",    "532955", -1.7,  568.95,  7635406,   "day", "13/08/2024 11:34:40"),
    @(10, "APOLLOTYRE", "Apollo Tyres Limited",                         "500877", -1.64, 485.4,   1443176,   "day", "13/08/2024 11:34:40"),
    @(11, "INDUSTOWER", "Indus Towers Ltd (Bharti Infratel)",           "534816", -2.55, 405.3,   7084711,   "day", "13/08/2024 11:34:40"),
    @(12, "BPCL",       "Bharat Petroleum Corporation Limited",         "500547", -3.51, 321.7,   9757429,   "day", "13/08/2024 11:34:40"),
    @(13, "BANDHANBNK", "Bandhan Bank Ltd",                             "541153", -2.26, 192.6,   14976179,  "day", "13/08/2024 11:34:40"),
    @(14, "IOC",        "Indian Oil Corporation Limited",               "530965", -2.98, 164.12,  16620160,  "day", "13/08/2024 11:34:40")
)

$startRow = 337
for ($idx = 0; $idx -lt $newRows.Count; $idx++) {
    $r = $startRow + $idx
    $data = $newRows[$idx]

    $ws.Cells.Item($r, 1).Value = $data[0]          # sr
    $ws.Cells.Item($r, 2).Value = $data[1]          # nsecode
    $ws.Cells.Item($r, 3).Value = $data[2]          # name
    $ws.Cells.Item($r, 4).Value = "'" + $data[3]    # bsecode (kept as text)
    $ws.Cells.Item($r, 5).Value = $data[4]          # per_chg
    $ws.Cells.Item($r, 6).Value = $data[5]          # close
    $ws.Cells.Item($r, 7).Value = $data[6]          # volume
    $ws.Cells.Item($r, 8).Value = $data[7]          # timeframe
    $ws.Cells.Item($r, 9).Value = $data[8]          # Date Time
}
